$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 43; existing rows 43-118 shift down to 44-119.
$ws.Rows.Item(43).Insert()

# Populate the new row 43 with the new weekly price entry.
$ws.Cells.Item(43, 1).Value = 8
$ws.Cells.Item(43, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(43, 3).Value = "Coquimbo"
$ws.Cells.Item(43, 4).Value = 44477
$ws.Cells.Item(43, 5).Value = 4
$ws.Cells.Item(43, 6).Value = 100112037
$ws.Cells.Item(43, 7).Value = "Cebollín"
$ws.Cells.Item(43, 8).Value = "Sin especificar"
$ws.Cells.Item(43, 9).Value = "Primera"
$ws.Cells.Item(43, 10).Value = 3200
$ws.Cells.Item(43, 11).Value = 900
$ws.Cells.Item(43, 12).Value = 1000
$ws.Cells.Item(43, 13).Value = 950
$ws.Cells.Item(43, 14).Value = "$/paquete 6 unidades"
$ws.Cells.Item(43, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(43, 16).Value = 158
$ws.Cells.Item(43, 17).Value = 6
$ws.Cells.Item(43, 18).Value = "Hortaliza"
